# "add jenkins build details as parameterize"
# The Registration sheet's "Regression" column (C) is switched from "YES"
# to "NO" for the bulk of the test-case rows (9-69), turning those rows
# off for regression runs (only the header block stays as-is).
# The active selection also moves from C58 up to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")
$ws.Activate()

$ws.Range("C9:C69").Value = "NO"

$ws.Range("C8").Select()
